# Apply the coin-price refresh captured in the commit:
# text/link/percentage cells are plain strings; price cells that look
# like plain numbers ("604.89", "1.00", ...) must be forced to stay text
# (NumberFormat "@") so Excel does not coerce them into numeric values,
# then restored to the default "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Plain text cells (coin names, links, formatted price strings, percentages)
$plainUpdates = [ordered]@{
    "D2" = "66.216.50"
    "E2" = "  -0.63%  "
    "D3" = "3.557.66"
    "E3" = "  +2.42%  "
    "E4" = "  +0.01%  "
    "E5" = "  +0.11%  "
    "E6" = "  -0.65%  "
    "D7" = "3.556.53"
    "E7" = "  +2.33%  "
    "E8" = "  +0.06%  "
    "E9" = "  +1.27%  "
    "E10" = "  +2.72%  "
    "E11" = "  -2.88%  "
    "E12" = "  -1.24%  "
    "D13" = "4.162.99"
    "E13" = "  +2.07%  "
    "E14" = "  -1.58%  "
    "E15" = "  -2.50%  "
    "D16" = "3.560.01"
    "E16" = "  +2.37%  "
    "D17" = "66.315.42"
    "E17" = "  -0.81%  "
    "E19" = "  +7.63%  "
    "E20" = "  -0.63%  "
    "E21" = "  -1.98%  "
    "E22" = "  -0.56%  "
    "E23" = "  +0.98%  "
    "E24" = "  -0.94%  "
    "D25" = "3.698.91"
    "E25" = "  +2.30%  "
    "E26" = "  -0.04%  "
    "E27" = "  +3.75%  "
    "E28" = "  +0.09%  "
    "E29" = "  -5.21%  "
    "E30" = "  -0.07%  "
    "E31" = "  -0.22%  "
    "E32" = "  -2.88%  "
    "E33" = "  -3.51%  "
    "B34" = "RenzoRestakedETH"
    "C34" = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
    "D34" = "3.553.61"
    "E34" = "  +2.15%  "
    "B35" = "EthereumClassic"
    "C35" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "E35" = "  +0.79%  "
    "E36" = "  -0.05%  "
    "E37" = "  -0.69%  "
    "B38" = "NEARProtocol"
    "C38" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "E38" = "  -1.66%  "
    "B39" = "Aptos"
    "C39" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "E39" = "  -0.78%  "
    "E40" = "  -0.13%  "
    "E41" = "  -1.05%  "
    "E42" = "  -3.53%  "
    "E43" = "  -0.37%  "
    "E44" = "  +0.77%  "
    "E45" = "  -3.95%  "
    "E46" = "  -1.53%  "
    "B47" = "InjectiveProtocol"
    "C47" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "E47" = "  -5.38%  "
    "B48" = "ONDO"
    "C48" = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
    "E48" = "  +1.02%  "
    "E49" = "  +0.45%  "
    "E50" = "  -1.89%  "
    "B51" = "SuiNetwork"
    "C51" = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
    "E51" = "  -3.30%  "
}
foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

# Price cells whose text looks like a plain number and must stay text
$numericTextUpdates = [ordered]@{
    "D5" = "604.89"
    "D6" = "144.60"
    "D10" = "8.08"
    "D11" = "0.136"
    "D14" = "0.0000207"
    "D15" = "30.18"
    "D19" = "11.44"
    "D20" = "6.20"
    "D21" = "14.88"
    "D22" = "427.80"
    "D23" = "0.605"
    "D24" = "78.74"
    "D27" = "0.0000120"
    "D28" = "8.05"
    "D29" = "9.18"
    "D30" = "2.48"
    "D35" = "25.42"
    "D38" = "5.63"
    "D39" = "7.81"
    "D40" = "1.00"
    "D41" = "172.67"
    "D42" = "0.0855"
    "D43" = "5.29"
    "D44" = "0.893"
    "D45" = "1.90"
    "D46" = "45.60"
    "D47" = "26.00"
    "D48" = "1.21"
    "D49" = "2.39"
    "D50" = "7.13"
    "D51" = "0.939"
}
foreach ($cellRef in $numericTextUpdates.Keys) {
    Set-TextCell $cellRef $numericTextUpdates[$cellRef]
}
